# update dataset to 16 april (adds 2020-04-16 / serial 43937 row to each sheet)

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Confirmed (sheet1): new row 41, B41 continues the running-total formula
# ---------------------------------------------------------------------------
$wsConfirmed = $wb.Worksheets.Item("Confirmed")
$wsConfirmed.Activate()

$wsConfirmed.Range("A40:C40").Copy()
$wsConfirmed.Range("A41:C41").PasteSpecial($xlPasteFormats)

$wsConfirmed.Range("A41").Value = 43937
$wsConfirmed.Range("B41").Formula = "=SUM(B40+C41)"
$wsConfirmed.Range("C41").Value = 341

# ---------------------------------------------------------------------------
# Recoverd (sheet2): new row 41
# ---------------------------------------------------------------------------
$wsRecoverd = $wb.Worksheets.Item("Recoverd")
$wsRecoverd.Activate()

$wsRecoverd.Range("A40:C40").Copy()
$wsRecoverd.Range("A41:C41").PasteSpecial($xlPasteFormats)

$wsRecoverd.Range("A41").Value = 43937
$wsRecoverd.Range("B41").Formula = "=SUM(B40+C41)"
$wsRecoverd.Range("C41").Value = 0

# ---------------------------------------------------------------------------
# Death (sheet3): new row 41
# ---------------------------------------------------------------------------
$wsDeath = $wb.Worksheets.Item("Death")
$wsDeath.Activate()

$wsDeath.Range("A40:C40").Copy()
$wsDeath.Range("A41:C41").PasteSpecial($xlPasteFormats)

$wsDeath.Range("A41").Value = 43937
$wsDeath.Range("B41").Formula = "=SUM(B40+C41)"
$wsDeath.Range("C41").Value = 10

# ---------------------------------------------------------------------------
# View state: selections / scroll position / active sheet for each sheet,
# matching what the workbook looked like after scrolling to the bottom of
# the (now one row longer) tables.
# ---------------------------------------------------------------------------
$wsRecoverd.Activate()
$wsRecoverd.Range("B40:B41").Select()

$wsDeath.Activate()
$wsDeath.Range("B40:B41").Select()

$wsConfirmed.Activate()
$wsConfirmed.Range("D44").Select()
